$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = 736482.01721401198
$ws.Range("D5").Value = 3750779.9454453602

$ws.Range("E4").Select()
